$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values (16/06/2025 and 23/06/2025 totals revised)
$ws.Cells.Item(12, 2).Value = 24877.27
$ws.Cells.Item(17, 2).Value = 69364.86

# Insert a new row for 24/06/2025 right after the existing 23/06/2025 row,
# shifting all subsequent rows down by one.
$ws.Rows.Item(18).Insert()
$ws.Cells.Item(18, 1).Value = 24
$ws.Cells.Item(18, 2).Value = 13051.24
$ws.Cells.Item(18, 3).Value = 6
$ws.Cells.Item(18, 4).Value = 2025
$ws.Cells.Item(18, 5).Value = "06/2025"
